$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header) ---
$ws.Range("A1").Value = "Projeto"
$ws.Range("B1").Value = "Defense in Depth"
$ws.Range("C1").Value = "Distrust Input"
$ws.Range("D1").Value = "Secure By Default"
$ws.Range("E1").Value = "Too Many Cooks"
$ws.Range("F1").Value = "Complex Inputs"
$ws.Range("G1").Value = "Fix Untested"
$ws.Range("H1").Value = "Code Refactors"
$ws.Range("I1").Value = "Serial Killer"
$ws.Range("J1").Value = "Lacked Test"
$ws.Range("K1").Value = "Frameworks are Optional"
$ws.Range("L1").Value = "Security By Obscurity"
$ws.Range("M1").Value = "You Ain't Gonna Need It"
$ws.Range("N1").Value = "Least Privilege"
$ws.Range("O1").Value = "Native Wrappers"

# Copy the header style from an existing styled header cell (A1) to the
# newly added header cells so H1/K1/N1/O1 match the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("H1:O1").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 2 (data) ---
$ws.Range("A2").Value = "Struts"
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 37
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 32
$ws.Range("F2").Value = 16
$ws.Range("G2").Value = 26
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 9
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 3
$ws.Range("M2").Value = 3
$ws.Range("N2").Value = 2
$ws.Range("O2").Value = 2
